# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the data table (row 136), pushing the
# existing records (rows 136-172) down by one row to rows 137-173, and populate
# the new row 136 with this week's observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (136:172) down by one to make room for the new record.
$ws.Rows("136:136").Insert()

# Populate the newly inserted row 136 with the new weekly record.
$ws.Cells.Item(136, 1).Value = 5
$ws.Cells.Item(136, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(136, 3).Value = "Maule"
$ws.Cells.Item(136, 4).Value = 44463
$ws.Cells.Item(136, 5).Value = 7
$ws.Cells.Item(136, 6).Value = 100112009
$ws.Cells.Item(136, 7).Value = "Acelga"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 400
$ws.Cells.Item(136, 11).Value = 2000
$ws.Cells.Item(136, 12).Value = 2000
$ws.Cells.Item(136, 13).Value = 2000
$ws.Cells.Item(136, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(136, 15).Value = "Región del Maule"
$ws.Cells.Item(136, 16).Value = 500
$ws.Cells.Item(136, 17).Value = 4
$ws.Cells.Item(136, 18).Value = "Hortaliza"
